$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-25 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-26 Sunday", 2) | Out-Null
$d.Content.Find.Execute("98÷6=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "38÷6=6, 2", 2) | Out-Null
$d.Content.Find.Execute("69÷7=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "62÷5=12, 2", 2) | Out-Null
$d.Content.Find.Execute("43÷2=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "74÷7=10, 4", 2) | Out-Null
$d.Content.Find.Execute("66÷8=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "61÷8=7, 5", 2) | Out-Null
$d.Content.Find.Execute("33÷7=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "33÷9=3, 6", 2) | Out-Null
$d.Content.Find.Execute("69÷9=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "44÷4=11, 0", 2) | Out-Null
$d.Content.Find.Execute("72÷9=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "40÷2=20, 0", 2) | Out-Null
$d.Content.Find.Execute("51÷4=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "40÷2=20, 0", 2) | Out-Null
$d.Content.Find.Execute("98÷4=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "51÷7=7, 2", 2) | Out-Null
$d.Content.Find.Execute("67÷3=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "56÷6=9, 2", 2) | Out-Null
$d.Content.Find.Execute("48÷5=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "97÷5=19, 2", 2) | Out-Null
$d.Content.Find.Execute("66÷7=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "12÷8=1, 4", 2) | Out-Null
$d.Content.Find.Execute("94÷3=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=5, 0", 2) | Out-Null
$d.Content.Find.Execute("85÷5=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=3, 7", 2) | Out-Null
$d.Content.Find.Execute("85÷3=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "36÷7=5, 1", 2) | Out-Null
$d.Content.Find.Execute("37÷5=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=31, 1", 2) | Out-Null
$d.Content.Find.Execute("99÷7=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "75÷9=8, 3", 2) | Out-Null
$d.Content.Find.Execute("83÷3=27, 2", $true, $false, $false, $false, $false, $true, 1, $false, "32÷6=5, 2", 2) | Out-Null
$d.Content.Find.Execute("99÷6=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "87÷4=21, 3", 2) | Out-Null
$d.Content.Find.Execute("26÷5=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "62÷9=6, 8", 2) | Out-Null
$d.Content.Find.Execute("96÷7=13, 5", $true, $false, $false, $false, $false, $true, 1, $false, "99÷7=14, 1", 2) | Out-Null
$d.Content.Find.Execute("50÷6=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "13÷3=4, 1", 2) | Out-Null
$d.Content.Find.Execute("10÷4=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "77÷7=11, 0", 2) | Out-Null
$d.Content.Find.Execute("52÷6=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "63÷7=9, 0", 2) | Out-Null
$d.Content.Find.Execute("72÷2=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "39÷7=5, 4", 2) | Out-Null
